# Refresh the crypto price/volume table to match the latest scrape.
# Prices in column D are stored as literal text (e.g. "18.30", "1.001")
# so a leading apostrophe is used to stop Excel from auto-converting
# them to numbers (which would silently drop trailing zeros / collapse
# the "x.y.zz" grouped-thousands strings used for the larger coins).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.476.58'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '''1.878.35'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''0.7143'
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('D6').Value = '''242.29'
$ws.Range('E6').Value = '  +1.63%  '
$ws.Range('D7').Value = '''1.001'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.07861'
$ws.Range('E8').Value = '  -1.76%  '
$ws.Range('D9').Value = '''0.3118'
$ws.Range('E9').Value = '  +2.89%  '
$ws.Range('D10').Value = '''25.22'
$ws.Range('E10').Value = '  +7.25%  '
$ws.Range('D11').Value = '''0.08268'
$ws.Range('D12').Value = '''0.7319'
$ws.Range('E12').Value = '  +3.42%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''5.281'
$ws.Range('E13').Value = '  +1.47%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '''1.844.00'
$ws.Range('E14').Value = '  -1.89%  '
$ws.Range('D15').Value = '''91.27'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '''29.431.38'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '''5.932'
$ws.Range('E17').Value = '  +1.57%  '
$ws.Range('D18').Value = '''247.70'
$ws.Range('E18').Value = '  +3.96%  '
$ws.Range('D19').Value = '''0.000007881'
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').Value = '''13.30'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').Value = '''0.9993'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('D22').Value = '''7.977'
$ws.Range('E22').Value = '  +6.44%  '
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '''0.1593'
$ws.Range('E24').Value = '  +10.54%  '
$ws.Range('D25').Value = '''163.89'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('D26').Value = '''9.017'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').Value = '''18.30'
$ws.Range('D28').Value = '''1.363'
$ws.Range('E28').Value = '  -4.44%  '
$ws.Range('D29').Value = '''1.496'
$ws.Range('E29').Value = '  +1.31%  '
$ws.Range('D30').Value = '''4.374'
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').Value = '''4.133'
$ws.Range('E31').Value = '  +2.51%  '
$ws.Range('D32').Value = '''0.05319'
$ws.Range('E32').Value = '  +2.33%  '
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('E34').Value = '  +3.38%  '
$ws.Range('D35').Value = '''0.7250'
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('D36').Value = '''2.680'
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '''0.01870'
$ws.Range('E37').Value = '  +0.85%  '
$ws.Range('D38').Value = '''1.263.94'
$ws.Range('E38').Value = '  +9.07%  '
$ws.Range('D39').Value = '''2.736'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').Value = '''0.9133'
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('D41').Value = '''74.00'
$ws.Range('E41').Value = '  +4.56%  '
$ws.Range('D42').Value = '''6.128'
$ws.Range('E42').Value = '  +1.91%  '
$ws.Range('D43').Value = '''1.001'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '''103.53'
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('D45').Value = '''0.5338'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('B46').Value = 'SynthetixNetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D46').Value = '''2.963'
$ws.Range('E46').Value = '  +13.98%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '''1.775'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').Value = '''0.4329'
$ws.Range('E49').Value = '  +1.33%  '
$ws.Range('D50').Value = '''9.260'
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('D51').Value = '''7.098'
$ws.Range('E51').Value = '  +1.51%  '
